{"js": "// Fix typos in the milestone paragraph.\n// Each entry is an unambiguous \"find\" substring (unique in the document)\n// paired with its replacement text.\nconst replacements = [\n  [\n    \"word and the letters change color\",\n    \"word and letters that change color\"\n  ],\n  [\n    \"like Duck Hunt, where birds move all about the screen\",\n    \"like Duck Hunt where birds move across the screen\"\n  ],\n  [\n    \"each with a word, they fall to the bottom if their word is completely typed and the game is lost of any mole reaches the top of the screen. Right now, we gave the groundwork\",\n    \"each with a word. The moles fall to the bottom if their word is completely typed and the game is lost if any mole reaches the top of the screen. Right now, we have the groundwork\"\n  ],\n  [\n    \"Having an on-screen keyboard is of high priority but should be achievable since we have experience doing that.\",\n    \"Having an on-screen keyboard is of high priority and should be easily achievable since we have experience doing similar tasks.\"\n  ],\n  [\n    \"a screen shot of\",\n    \"a screen-shot of\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix typos in the milestone paragraph.\n\nfunction Replace-DocText {\n    param(\n        $Document,\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $find = $Document.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n\n    # wdReplaceOne = 1 replaces a single match.\n    $find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 1) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\nReplace-DocText $d \"word and the letters change color\" \"word and letters that change color\"\nReplace-DocText $d \"like Duck Hunt, where birds move all about the screen\" \"like Duck Hunt where birds move across the screen\"\nReplace-DocText $d \"each with a word, they fall to the bottom if their word is completely typed and the game is lost of any mole reaches the top of the screen. Right now, we gave the groundwork\" \"each with a word. The moles fall to the bottom if their word is completely typed and the game is lost if any mole reaches the top of the screen. Right now, we have the groundwork\"\nReplace-DocText $d \"Having an on-screen keyboard is of high priority but should be achievable since we have experience doing that.\" \"Having an on-screen keyboard is of high priority and should be easily achievable since we have experience doing similar tasks.\"\nReplace-DocText $d \"a screen shot of\" \"a screen-shot of\"\n"}
